$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.9998630581663609
$ws.Range("E3").Value = 0.9998630581663609

$ws.Range("D4").Value = 0.000002015437435285352
$ws.Range("E4").Value = 0.000002015437435285352

$ws.Range("D5").Value = 0.0000000003353121724228589
$ws.Range("E5").Value = 0.0000000003353121724228589

$ws.Range("D6").Value = 0.0000000000000000004100039910519209
$ws.Range("E6").Value = 0.0000000000000000004100039910519209

$ws.Range("F7").Value = 12.50727081298828
